$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F9").Value = 2173
$ws.Range("F10").Value = 623
$ws.Range("F13").Value = 1065
$ws.Range("F15").Value = 2198
$ws.Range("F16").Value = 654
$ws.Range("F17").Value = 12622
$ws.Range("F18").Value = 1236
$ws.Range("F19").Value = 9
$ws.Range("F20").Value = 556
$ws.Range("F21").Value = 127
$ws.Range("F22").Value = 22
$ws.Range("F25").Value = 264

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F21").Value = 4

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5698
$ws.Range("F4").Value = 466

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5698
$ws.Range("F5").Value = 466
$ws.Range("F13").Value = 2173
$ws.Range("F14").Value = 623
$ws.Range("F19").Value = 1065
$ws.Range("F24").Value = 2198
$ws.Range("F25").Value = 654
$ws.Range("F28").Value = 1236
$ws.Range("F29").Value = 9
$ws.Range("F30").Value = 556
$ws.Range("F31").Value = 127
$ws.Range("F32").Value = 22
$ws.Range("F38").Value = 264
$ws.Range("F45").Value = 4
